# Apply the "Contacts data" sample-data update to the workbook.

$wb = $excel.ActiveWorkbook

# 1. Rename the "data" worksheet to "Contacts data"
$ws = $wb.Worksheets.Item("data")
$ws.Name = "Contacts data"

# 2. Clear the previous sample contents and write the new table.
$ws.Cells.Clear()

$headers = @("Id", "FirstName", "LastName", "Email", "Phone", "AccountId")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$rows = @(
    @("0038N00000D0REQQA3", "Jon",    "Amos",    "info@salesforce.com",     "(905) 555-1212", "0018N00000EyKhEQAV"),
    @("0038N00000D0RERQA3", "John",   "Smith",   "john@servicecloud.trial", "(212) 555-5555", "0018N00000EyKhFQAV"),
    @("0038N00000D0RESQA3", "Geoff",  "Minor",   "info@salesforce.com",     "(415) 555-1212", "0018N00000EyKhEQAV"),
    @("0038N00000D0RETQA3", "Carole", "White",   "info@salesforce.com",     "(415) 555-1212", "0018N00000EyKhEQAV"),
    @("0038N00000D0REUQA3", "Edward", "Stamos",  "info@salesforce.com",     "(212) 555-5555", "0018N00000EyKhFQAV"),
    @("0038N00000D0REVQA3", "Howard", "Jones",   "info@salesforce.com",     "(212) 555-5555", "0018N00000EyKhFQAV"),
    @("0038N00000D0REWQA3", "Leanne", "Tomlin",  "info@salesforce.com",     "(212) 555-5555", "0018N00000EyKhFQAV"),
    @("0038N00000D0REXQA3", "Marc",   "Benioff", "info@salesforce.com",     "(415) 901-7000", "0018N00000EyKhGQAV"),
    @("0038N00000D0REYQA3", "Mark",   "Land",    "your.email@test.com",     "(415) 555-5555", "0018N00000EyKhFQAV")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
